$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (travis)
$ws.Range("B2").Value = 212
$ws.Range("C2").Value = 587
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 94
$ws.Range("F2").Value = 2.77
$ws.Range("G2").Value = 2

# Row 3 (github actions)
$ws.Range("B3").Value = 495
$ws.Range("C3").Value = 3508
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 312
$ws.Range("F3").Value = 7.09
$ws.Range("G3").Value = 4

# Row 4 (gitlab ci)
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 44
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 4.89
$ws.Range("G4").Value = 4
